$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the last_update timestamp for row 5 (bevnat_info)
$ws.Range("E5").Value = 1706219962

# Move/update the active selection to E9
$ws.Range("E9").Select()
